# Apply targeted "statement_sub_section" corrections on the "cbs_38" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cbs_38")

# E8:E11 -> noncurrent (was current)
$ws.Range("E8").Value = "noncurrent"
$ws.Range("E9").Value = "noncurrent"
$ws.Range("E10").Value = "noncurrent"
$ws.Range("E11").Value = "noncurrent"

# Row 13: D13 -> equity_liabilities (was assets), E13 -> current (was noncurrent)
$ws.Range("D13").Value = "equity_liabilities"
$ws.Range("E13").Value = "current"

# E21:E23 -> noncurrent (was current)
$ws.Range("E21").Value = "noncurrent"
$ws.Range("E22").Value = "noncurrent"
$ws.Range("E23").Value = "noncurrent"

# E25:E26 -> equity (was noncurrent)
$ws.Range("E25").Value = "equity"
$ws.Range("E26").Value = "equity"
